# Update BAMM codes and related descriptions.
#
# This script updates several "R Packages" bullet items in the document to
# insert explicit version numbers for BAMMtools, phytools, geiger, ape, and
# PhylogeneticEM, and fixes the "PhylogeneticEN" -> "PhylogeneticEM" typo.

$d = $word.ActiveDocument

function Replace-InParagraphRange([object]$range, [string]$findText, [string]$replaceText) {
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $found = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    return $found
}

$paras = $d.Paragraphs
$n = $paras.Count

for ($i = 1; $i -le $n; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text

    if ($t.StartsWith("ape:")) {
        # "ape: Functions used: ..." -> "ape v5.7.1: Functions used: ..."
        Replace-InParagraphRange $p.Range ": Functions used: " " v5.7.1: Functions used: " | Out-Null
    }
    elseif ($t.StartsWith("BAMMtools:")) {
        # "BAMMtools: Functions used: ..." -> "BAMMtools v2.1.9: Functions used: ..."
        Replace-InParagraphRange $p.Range ": Functions used: " " v2.1.9: Functions used: " | Out-Null
    }
    elseif ($t.StartsWith("phytools:")) {
        # "phytools: Functions used: ..." -> "phytools v1.5.1: Functions used: ..."
        Replace-InParagraphRange $p.Range ": Functions used: " " v1.5.1: Functions used: " | Out-Null
    }
    elseif ($t.StartsWith("geiger:")) {
        # "geiger: Functions used: ..." -> "geiger v2.0.9: Functions used: ..."
        Replace-InParagraphRange $p.Range ": Functions used: " " v2.0.9: Functions used: " | Out-Null
    }
    elseif ($t.StartsWith("PhylogeneticEN:")) {
        # "PhylogeneticEN: Function used: ..." -> "PhylogeneticEM v1.6.0: Function used: ..."
        Replace-InParagraphRange $p.Range "PhylogeneticEN" "PhylogeneticEM" | Out-Null
        Replace-InParagraphRange $p.Range ": Function used: " " v1.6.0: Function used: " | Out-Null
    }
    elseif ($t.StartsWith("OUwie:")) {
        # "OUwie: Used for statistical modeling." -> "OUwie v2.10: Used for statistical modeling."
        Replace-InParagraphRange $p.Range ": Used for statistical modeling." " v2.10: Used for statistical modeling." | Out-Null
    }
}
